# Update the F-column (ticket/view count) figures on the "展览" and
# "全部类型" worksheets to reflect the latest generated snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 4,5,6,7,9,10,11,13
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 275
$wsExhibit.Range("F5").Value = 154
$wsExhibit.Range("F6").Value = 71
$wsExhibit.Range("F7").Value = 272
$wsExhibit.Range("F9").Value = 1998
$wsExhibit.Range("F10").Value = 352
$wsExhibit.Range("F11").Value = 4721
$wsExhibit.Range("F13").Value = 332

# Sheet "全部类型" (All types) - same events, rows 6,7,8,9,13,14,15,17
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 275
$wsAll.Range("F7").Value = 154
$wsAll.Range("F8").Value = 71
$wsAll.Range("F9").Value = 272
$wsAll.Range("F13").Value = 1998
$wsAll.Range("F14").Value = 352
$wsAll.Range("F15").Value = 4721
$wsAll.Range("F17").Value = 332
